# Update the dSF column (F) values to reflect a repull/recalculation of the data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F4").Value = -9
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = 2
